$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - only Frequency changes
$ws.Range("B2").Value = 101

# Row 3
$ws.Range("A3").Value = "bayram, atatürk, yıl, gazi, ecdat"
$ws.Range("B3").Value = 40
$ws.Range("C3").Value = "bayram, atatürk, türk, kutlu, mustafa"

# Row 4
$ws.Range("A4").Value = "millet, sandık, seçim, oy, mayıs"
$ws.Range("B4").Value = 35
$ws.Range("C4").Value = "millet, oy, buluş, sandık, seçim"

# Row 5
$ws.Range("A5").Value = "milyar, yatırım, lira, kamu, dolar"
$ws.Range("B5").Value = 31
$ws.Range("C5").Value = "milyar, yatırım, lira, dolar, kamu"

# Row 6
$ws.Range("A6").Value = "türkiye, yüzyıl, türk, milliyetçilik, millet"
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = "türk, türkiye, milliyetçi, yüzyıl, millet"

# Row 7 - Frequency unchanged (26)
$ws.Range("A7").Value = "muhteşem, bil, şanlıurfa, geleneksel, van"
$ws.Range("C7").Value = "muhteşem, bil, şanlıurfa, van, geleneksel"

# Row 8
$ws.Range("A8").Value = "teşekkür, muhteşem, istanbul, büyük, kardeş"
$ws.Range("B8").Value = 25
$ws.Range("C8").Value = "teşekkür, başkan, dernek, muhteşem, misafirperverlikleri"

# Row 9 - Frequency unchanged (24)
$ws.Range("A9").Value = "nükleer, santral, enerji, üretim, gaz"
$ws.Range("C9").Value = "nükleer, santral, enerji, üretim, gaz"

# Row 10 - Frequency unchanged (23)
$ws.Range("A10").Value = "şehit, rahmet, atatürk, emanet, havaliman"
$ws.Range("C10").Value = "şehit, rahmet, atatürk, an, dönüm"

# Row 11
$ws.Range("A11").Value = "canlı, yayın, ortak, bölüm, tv"
$ws.Range("B11").Value = 18
$ws.Range("C11").Value = "yayın, program, canlı, konuk, sun"

# Row 12 - Frequency unchanged (15)
$ws.Range("A12").Value = "gök, iş, gülsün, allah, sevdalı"
$ws.Range("C12").Value = "iyi, sanatçı, allah, iş, çık"

# Row 13
$ws.Range("A13").Value = "ırak, terör, kardeş, örgüt, karşı"
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = "ırak, terör, örgüt, karşı, kardeş"
